# Apply the updated crypto price/volume figures scraped on
# Sat Sep 16 21:39:33 UTC 2023, matching the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.737.76"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.646.75"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.26"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.876.48"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "1.662.47"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.23"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "26.746.18"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.14"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("E21").Value = "  +14.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.28"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.46"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.56"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +3.94%  "
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "1.280.37"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.540"
$ws.Range("E38").Value = "  +5.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.830"
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.816"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("D44").Value = "1.787.97"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.98"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.89"
$ws.Range("E46").Value = "  +8.90%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0987"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("E51").Value = "  +2.11%  "
